# Commit: "Add FsSpreadsheet tests and update testFile to include Time in DateTime"
#
# 1) Give the "DateTime" column in the WithTable sheet a couple of
#    timestamps that actually carry a time component (previously every
#    value was midnight), using a new custom number format that also
#    shows the time (d/m/yy h:mm).
# 2) Widen column C on WithTable to fit the new format.
# 3) Update the saved cursor/selection on the WithTable and Tableless
#    sheets.
# 4) Leave WithTable (the first sheet) as the active tab again.

$wb = $excel.ActiveWorkbook

$wsWithTable  = $wb.Worksheets.Item("WithTable")
$wsTableless  = $wb.Worksheets.Item("Tableless")

# --- DateTime column: add a time-of-day component to two rows ---------
# Row 3 (C3): 2023-10-15 00:00  ->  2023-10-15 18:00
$wsWithTable.Range("C3").NumberFormat = "d/m/yy h:mm;@"
$wsWithTable.Range("C3").Value = 45214.75

# Row 4 (C4): 2023-10-16 00:00  ->  2023-10-16 20:00
$wsWithTable.Range("C4").NumberFormat = "d/m/yy h:mm;@"
$wsWithTable.Range("C4").Value = 45215.83333333333

# --- Column width: widen column C (DateTime) so the new format fits ---
$wsWithTable.Columns.Item(3).ColumnWidth = 25.022135416666668

# --- Selection / active sheet bookkeeping ------------------------------
# Touch the Tableless sheet first so it is no longer the last-focused
# sheet, then finish on WithTable so it ends up active again (tab 0).
[void]$wsTableless.Range("F1").Select()
[void]$wsWithTable.Range("E11").Select()
